$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 250.375
$ws.Range("I4").Value = 250.375
$ws.Range("K4").Value = 250.375
$ws.Range("M4").Value = -136.375
$ws.Range("H17").Value = 376.7143
$ws.Range("J17").Value = 376.7143
$ws.Range("L17").Value = 1130.1429
$ws.Range("N17").Value = -1466.1429
$ws.Range("H18").Value = 3596.6
$ws.Range("I18").Value = 3370.75
$ws.Range("K18").Value = 3370.75
$ws.Range("M18").Value = -3086.75
$ws.Range("H19").Value = 1376.3
$ws.Range("I19").Value = 1319.75
$ws.Range("J19").Value = 1414
$ws.Range("K19").Value = 1319.75
$ws.Range("L19").Value = 1414
$ws.Range("M19").Value = -1144.75
$ws.Range("N19").Value = -1764
$ws.Range("H28").Value = 1074.9546
$ws.Range("I28").Value = 455.6316
$ws.Range("K28").Value = 455.6316
$ws.Range("M28").Value = 29.36840000000001
$ws.Range("H33").Value = 493.08334
$ws.Range("I33").Value = 363.26315
$ws.Range("J33").Value = 986.4
$ws.Range("K33").Value = 363.26315
$ws.Range("L33").Value = 986.4
$ws.Range("M33").Value = -134.26315
$ws.Range("N33").Value = -1444.4
$ws.Range("H40").Value = 3399.8
$ws.Range("J40").Value = 4000
$ws.Range("L40").Value = 4000
$ws.Range("N40").Value = -4350
$ws.Range("H48").Value = 1000
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H51").Value = 7323.8667
$ws.Range("I51").Value = 5899
$ws.Range("J51").Value = 7543.077
$ws.Range("K51").Value = 5899
$ws.Range("L51").Value = 7543.077
$ws.Range("M51").Value = -5415
$ws.Range("N51").Value = -8511.077000000001
$ws.Range("H56").Value = 1000
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H64").Value = 6660.6665
$ws.Range("I64").Value = 4992.5
$ws.Range("K64").Value = 4992.5
$ws.Range("M64").Value = -4744.5
$ws.Range("H67").Value = 6660.6665
$ws.Range("I67").Value = 4992.5
$ws.Range("K67").Value = 4992.5
$ws.Range("M67").Value = -4134.5
$ws.Range("H86").Value = 1478.9584
$ws.Range("I86").Value = 1823.1111
$ws.Range("J86").Value = 1272.4667
$ws.Range("K86").Value = 1823.1111
$ws.Range("L86").Value = 1272.4667
$ws.Range("M86").Value = -700.1111000000001
$ws.Range("N86").Value = -3518.4667
$ws.Range("H88").Value = 1001601.5
$ws.Range("J88").Value = 3204
$ws.Range("L88").Value = 3204
$ws.Range("N88").Value = -4016
$ws.Range("H89").Value = 1478.9584
$ws.Range("I89").Value = 1823.1111
$ws.Range("J89").Value = 1272.4667
$ws.Range("K89").Value = 9115.5555
$ws.Range("L89").Value = 6362.3335
$ws.Range("M89").Value = -3499.5555
$ws.Range("N89").Value = -17594.3335
$ws.Range("H91").Value = 1001601.5
$ws.Range("J91").Value = 3204
$ws.Range("L91").Value = 3204
$ws.Range("N91").Value = -6012
$ws.Range("H107").Value = 1356.3334
$ws.Range("I107").Value = 1132.8334
$ws.Range("K107").Value = 1132.8334
$ws.Range("M107").Value = 787.1666
$ws.Range("H126").Value = 67490.336
$ws.Range("J126").Value = 67490.336
$ws.Range("L126").Value = 67490.336
$ws.Range("N126").Value = -77370.336
$ws.Range("H132").Value = 3795.7778
$ws.Range("I132").Value = 2881.1428
$ws.Range("K132").Value = 8643.428400000001
$ws.Range("M132").Value = -6113.428400000001
$ws.Range("H138").Value = 2190.38
$ws.Range("I138").Value = 1974.9412
$ws.Range("J138").Value = 2301.3635
$ws.Range("K138").Value = 5924.8236
$ws.Range("L138").Value = 6904.0905
$ws.Range("M138").Value = -784.8235999999997
$ws.Range("N138").Value = -17184.0905
$ws.Range("H141").Value = 5462.8887
$ws.Range("I141").Value = 3361
$ws.Range("J141").Value = 9666.666999999999
$ws.Range("K141").Value = 10083
$ws.Range("L141").Value = 29000.001
$ws.Range("M141").Value = -4903
$ws.Range("N141").Value = -39360.001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 550
$ws.Range("I4").Value = 550
$ws.Range("K4").Value = 550
$ws.Range("M4").Value = -434
$ws.Range("H5").Value = 1317.5454
$ws.Range("I5").Value = 999
$ws.Range("J5").Value = 2167
$ws.Range("K5").Value = 999
$ws.Range("L5").Value = 2167
$ws.Range("M5").Value = -887
$ws.Range("N5").Value = -2391
$ws.Range("H34").Value = 8364829
$ws.Range("I34").Value = 14312570
$ws.Range("J34").Value = 37990.8
$ws.Range("K34").Value = 14312570
$ws.Range("L34").Value = 37990.8
$ws.Range("M34").Value = -14312299
$ws.Range("N34").Value = -38532.8
$ws.Range("H45").Value = 1529.9375
$ws.Range("I45").Value = 1333.2
$ws.Range("K45").Value = 1333.2
$ws.Range("M45").Value = -956.2
$ws.Range("H61").Value = 4434.05
$ws.Range("I61").Value = 4315.722
$ws.Range("J61").Value = 5499
$ws.Range("K61").Value = 4315.722
$ws.Range("L61").Value = 5499
$ws.Range("M61").Value = -4103.722
$ws.Range("N61").Value = -5923
$ws.Range("H74").Value = 1762.8518
$ws.Range("I74").Value = 1448
$ws.Range("K74").Value = 1448
$ws.Range("M74").Value = -574
$ws.Range("H77").Value = 1762.8518
$ws.Range("I77").Value = 1448
$ws.Range("K77").Value = 7240
$ws.Range("M77").Value = -2872
$ws.Range("H97").Value = 2421.6
$ws.Range("I97").Value = 1784.5834
$ws.Range("J97").Value = 4969.6665
$ws.Range("K97").Value = 1784.5834
$ws.Range("L97").Value = 4969.6665
$ws.Range("M97").Value = -1288.5834
$ws.Range("N97").Value = -5961.6665
$ws.Range("H102").Value = 2114.5833
$ws.Range("I102").Value = 2114.5833
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2114.5833
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -492.5832999999998
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 4673.1
$ws.Range("I122").Value = 3432.75
$ws.Range("K122").Value = 10298.25
$ws.Range("M122").Value = -7848.25
$ws.Range("H136").Value = 4434.05
$ws.Range("I136").Value = 4315.722
$ws.Range("J136").Value = 5499
$ws.Range("K136").Value = 12947.166
$ws.Range("L136").Value = 16497
$ws.Range("M136").Value = -10397.166
$ws.Range("N136").Value = -21597

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1317.5454
$ws.Range("I4").Value = 999
$ws.Range("J4").Value = 2167
$ws.Range("K4").Value = 999
$ws.Range("L4").Value = 2167
$ws.Range("M4").Value = -884
$ws.Range("N4").Value = -2397
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H94").Value = 4863.6665
$ws.Range("I94").Value = 1070.875
$ws.Range("J94").Value = 12449.25
$ws.Range("K94").Value = 1070.875
$ws.Range("L94").Value = 12449.25
$ws.Range("M94").Value = -619.875
$ws.Range("N94").Value = -13351.25
$ws.Range("H97").Value = 24891.8
$ws.Range("I97").Value = 5997
$ws.Range("J97").Value = 100471
$ws.Range("K97").Value = 5997
$ws.Range("L97").Value = 100471
$ws.Range("M97").Value = -5006
$ws.Range("N97").Value = -102453
$ws.Range("H107").Value = 1093.4762
$ws.Range("J107").Value = 1090.7142
$ws.Range("L107").Value = 1090.7142
$ws.Range("N107").Value = -4930.7142
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H141").Value = 129998.5
$ws.Range("J141").Value = 129998.5
$ws.Range("L141").Value = 129998.5
$ws.Range("N141").Value = -140358.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2585.2856
$ws.Range("I16").Value = 2691.1667
$ws.Range("K16").Value = 2691.1667
$ws.Range("M16").Value = -2404.1667
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H31").Value = 2925.2727
$ws.Range("I31").Value = 1872.1428
$ws.Range("K31").Value = 1872.1428
$ws.Range("M31").Value = -1577.1428
$ws.Range("H34").Value = 2925.2727
$ws.Range("I34").Value = 1872.1428
$ws.Range("K34").Value = 1872.1428
$ws.Range("M34").Value = -1670.1428
$ws.Range("H62").Value = 3398
$ws.Range("I62").Value = 3398
$ws.Range("K62").Value = 3398
$ws.Range("M62").Value = -2774
$ws.Range("H65").Value = 3398
$ws.Range("I65").Value = 3398
$ws.Range("K65").Value = 16990
$ws.Range("M65").Value = -13870
$ws.Range("H99").Value = 26233.875
$ws.Range("J99").Value = 50326.91
$ws.Range("L99").Value = 50326.91
$ws.Range("N99").Value = -53322.91
$ws.Range("H105").Value = 1882.9565
$ws.Range("I105").Value = 1045.4445
$ws.Range("K105").Value = 1045.4445
$ws.Range("M105").Value = 701.5554999999999
$ws.Range("H113").Value = 2585.2856
$ws.Range("I113").Value = 2691.1667
$ws.Range("K113").Value = 2691.1667
$ws.Range("M113").Value = -521.1667000000002
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 26233.875
$ws.Range("J126").Value = 50326.91
$ws.Range("L126").Value = 150980.73
$ws.Range("N126").Value = -155920.73
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 52877.5
$ws.Range("H132").Value = 1689.3077
$ws.Range("I132").Value = 1714.6364
$ws.Range("K132").Value = 5143.9092
$ws.Range("M132").Value = -2613.9092
$ws.Range("H134").Value = 3175.2632
$ws.Range("I134").Value = 1850.0769
$ws.Range("K134").Value = 5550.2307
$ws.Range("M134").Value = -3015.2307

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 192
$ws.Range("J46").Value = 192
$ws.Range("L46").Value = 576
$ws.Range("N46").Value = -758
$ws.Range("H51").Value = 16165.667
$ws.Range("I51").Value = 16165.667
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 48497.001
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -48037.001
$ws.Range("N51").ClearContents()
$ws.Range("H80").Value = 2519000
$ws.Range("I80").Value = 38000
$ws.Range("K80").Value = 114000
$ws.Range("M80").Value = -113064
$ws.Range("H83").Value = 2519000
$ws.Range("I83").Value = 38000
$ws.Range("K83").Value = 342000
$ws.Range("M83").Value = -337320
$ws.Range("H120").Value = 2366.3333
$ws.Range("I120").Value = 2366.3333
$ws.Range("K120").Value = 7098.999899999999
$ws.Range("M120").Value = -2260.999899999999
$ws.Range("H131").Value = 20644.69
$ws.Range("I131").Value = 112248.5
$ws.Range("J131").Value = 1560.5625
$ws.Range("K131").Value = 336745.5
$ws.Range("L131").Value = 4681.6875
$ws.Range("M131").Value = -331705.5
$ws.Range("N131").Value = -14761.6875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 69028.39999999999
$ws.Range("J64").Value = 69028.39999999999
$ws.Range("L64").Value = 69028.39999999999
$ws.Range("N64").Value = -69524.39999999999
$ws.Range("H67").Value = 69028.39999999999
$ws.Range("J67").Value = 69028.39999999999
$ws.Range("L67").Value = 69028.39999999999
$ws.Range("N67").Value = -70744.39999999999
$ws.Range("H102").Value = 3327.087
$ws.Range("I102").Value = 2812.25
$ws.Range("J102").Value = 3888.7273
$ws.Range("K102").Value = 2812.25
$ws.Range("L102").Value = 3888.7273
$ws.Range("M102").Value = -1190.25
$ws.Range("N102").Value = -7132.7273
$ws.Range("H122").Value = 1411.091
$ws.Range("I122").Value = 1172.2
$ws.Range("K122").Value = 3516.6
$ws.Range("M122").Value = -1066.6
$ws.Range("H123").Value = 44415.793
$ws.Range("J123").Value = 49486.31
$ws.Range("L123").Value = 49486.31
$ws.Range("N123").Value = -54386.31
$ws.Range("H126").Value = 4569.25
$ws.Range("I126").Value = 3899
$ws.Range("J126").Value = 4792.6665
$ws.Range("K126").Value = 11697
$ws.Range("L126").Value = 14377.9995
$ws.Range("M126").Value = -9227
$ws.Range("N126").Value = -19317.9995
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4880.077
$ws.Range("I40").Value = 3064
$ws.Range("K40").Value = 3064
$ws.Range("M40").Value = -2928
$ws.Range("H54").Value = 34999.332
$ws.Range("I54").Value = 4999
$ws.Range("J54").Value = 49999.5
$ws.Range("K54").Value = 4999
$ws.Range("L54").Value = 49999.5
$ws.Range("M54").Value = -4355
$ws.Range("N54").Value = -51287.5
$ws.Range("H56").Value = 58399.6
$ws.Range("I56").Value = 54999
$ws.Range("K56").Value = 54999
$ws.Range("M56").Value = -54308
$ws.Range("H68").Value = 2635.077
$ws.Range("I68").Value = 2673
$ws.Range("K68").Value = 2673
$ws.Range("M68").Value = -1924
$ws.Range("H71").Value = 2635.077
$ws.Range("I71").Value = 2673
$ws.Range("K71").Value = 13365
$ws.Range("M71").Value = -9621
$ws.Range("H109").Value = 67332.336
$ws.Range("J109").Value = 67332.336
$ws.Range("L109").Value = 67332.336
$ws.Range("N109").Value = -70106.336
$ws.Range("H132").Value = 3617.25
$ws.Range("I132").Value = 3323.4285
$ws.Range("J132").Value = 4028.6
$ws.Range("K132").Value = 9970.2855
$ws.Range("L132").Value = 12085.8
$ws.Range("M132").Value = -7440.2855
$ws.Range("N132").Value = -17145.8
$ws.Range("H136").Value = 3151.0952
$ws.Range("I136").Value = 2953.75
$ws.Range("J136").Value = 3782.6
$ws.Range("K136").Value = 8861.25
$ws.Range("L136").Value = 11347.8
$ws.Range("M136").Value = -6311.25
$ws.Range("N136").Value = -16447.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 47499.5
$ws.Range("I40").Value = 47499.5
$ws.Range("K40").Value = 47499.5
$ws.Range("M40").Value = -47350.5
$ws.Range("H44").Value = 5000
$ws.Range("J44").Value = 5000
$ws.Range("L44").Value = 5000
$ws.Range("N44").Value = -6108
$ws.Range("H51").Value = 5535
$ws.Range("I51").Value = 5535
$ws.Range("K51").Value = 5535
$ws.Range("M51").Value = -5025
$ws.Range("H62").Value = 6192.222
$ws.Range("I62").Value = 3532.8572
$ws.Range("K62").Value = 3532.8572
$ws.Range("M62").Value = -2908.8572
$ws.Range("H65").Value = 6192.222
$ws.Range("I65").Value = 3532.8572
$ws.Range("K65").Value = 17664.286
$ws.Range("M65").Value = -14544.286
$ws.Range("H69").Value = 51710
$ws.Range("J69").Value = 51710
$ws.Range("L69").Value = 51710
$ws.Range("N69").Value = -53208
$ws.Range("H70").Value = 20818.092
$ws.Range("I70").Value = 16999.5
$ws.Range("K70").Value = 16999.5
$ws.Range("M70").Value = -16684.5
$ws.Range("H72").Value = 51710
$ws.Range("J72").Value = 51710
$ws.Range("L72").Value = 155130
$ws.Range("N72").Value = -162618
$ws.Range("H73").Value = 20818.092
$ws.Range("I73").Value = 16999.5
$ws.Range("K73").Value = 16999.5
$ws.Range("M73").Value = -15907.5
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H109").Value = 79949.25
$ws.Range("J109").Value = 79949.25
$ws.Range("L109").Value = 79949.25
$ws.Range("N109").Value = -82723.25
$ws.Range("H122").Value = 5125.486
$ws.Range("I122").Value = 6410.2104
$ws.Range("J122").Value = 3599.875
$ws.Range("K122").Value = 19230.6312
$ws.Range("L122").Value = 10799.625
$ws.Range("M122").Value = -16780.6312
$ws.Range("N122").Value = -15699.625
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 4468.4736
$ws.Range("I132").Value = 4438.9443
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 13316.8329
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -10786.8329
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 3539.8823
$ws.Range("I136").Value = 1298.5
$ws.Range("K136").Value = 3895.5
$ws.Range("M136").Value = -1345.5
